$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows (old rows 5-7 duplicated the FAPs sending-cluster
# data that's been folded into rows 2-4 with refreshed TPM-derived values).
$ws.Rows("5:7").Delete()

# Refresh rows 2-4 with the recomputed TPM-based values (new sending cluster
# "FAPs" / ligand "Wnt1" / receptor "Lrp6", recalculated specificity metrics).
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt1"
$ws.Range("C2").Value = "Lrp6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3532066666666667
$ws.Range("H2").Value = 1.05962
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 12.40685866666667
$ws.Range("N2").Value = 37.220576
$ws.Range("O2").Value = 0.1720325859617629
$ws.Range("P2").Value = 0.1720325859617629
$ws.Range("Q2").Value = 4.382185193457778
$ws.Range("R2").Value = 39.43966674112
$ws.Range("S2").Value = 0.1720325859617629
$ws.Range("T2").Value = 0.1720325859617629

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt1"
$ws.Range("C3").Value = "Lrp6"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3532066666666667
$ws.Range("H3").Value = 1.05962
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 43.97212233333332
$ws.Range("N3").Value = 131.916367
$ws.Range("O3").Value = 0.6097142007069145
$ws.Range("P3").Value = 0.6097142007069145
$ws.Range("Q3").Value = 15.53124675561555
$ws.Range("R3").Value = 139.78122080054
$ws.Range("S3").Value = 0.6097142007069145
$ws.Range("T3").Value = 0.6097142007069145

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt1"
$ws.Range("C4").Value = "Lrp6"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3532066666666667
$ws.Range("H4").Value = 1.05962
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 15.740255
$ws.Range("N4").Value = 47.220765
$ws.Range("O4").Value = 0.2182532133313226
$ws.Range("P4").Value = 0.2182532133313226
$ws.Range("Q4").Value = 5.559563001033333
$ws.Range("R4").Value = 50.0360670093
$ws.Range("S4").Value = 0.2182532133313226
$ws.Range("T4").Value = 0.2182532133313226

